# ajout ssr, had et psy
# Add one new row (46) at the bottom of the format description sheet:
#   B46 = "c"    (type)
#   D46 = 182    (position)
#   F46 = "ZAD"  (nom) -> becomes a new shared string entry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B46").Value = "c"
$ws.Range("D46").Value = 182
$ws.Range("F46").Value = "ZAD"

# Reflect the author's view state: scrolled down with D46 selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 8
$ws.Range("D46").Select()
